# Adds support for "MEF" (Managed Extensibility Framework) as a new row
# in the Tabelle1 benchmark table, between "LinFu" (row 8) and "Mugen"
# (row 9), shifting everything below down by one row, and updates the
# three 3-D bar charts' source ranges to include it. Also left-aligns
# column A (new cellXf + applied across A1:A19).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- 1. Insert a new row 9 (pushes old row 9 "Mugen" -> row 10, etc.) ---
$ws.Rows.Item(9).Insert()

# --- 2. Populate the new "MEF" row ---
$ws.Range("A9").Value = "MEF"
$ws.Range("B9").Value = 2838
$ws.Range("C9").Value = 9931
$ws.Range("D9").Value = 27133

# --- 3. Left-align all of column A's data (A1:A19) ---
$ws.Range("A1:A19").HorizontalAlignment = -4131
$ws.Range("A19").Select()

# --- 4. Re-point the three embedded charts at the grown (and shifted)
#        non-contiguous ranges, skipping the "Spring.NET" outlier row
#        (now row 15 instead of row 14). ---
$catFormula = "(Tabelle1!`$A`$2:`$A`$14,Tabelle1!`$A`$16:`$A`$19)"

$chart1 = $ws.ChartObjects(1).Chart
$chart1.SeriesCollection(1).Formula = "=SERIES(,(Tabelle1!`$A`$2:`$A`$14,Tabelle1!`$A`$16:`$A`$19),(Tabelle1!`$B`$2:`$B`$14,Tabelle1!`$B`$16:`$B`$19),1)"

$chart2 = $ws.ChartObjects(2).Chart
$chart2.SeriesCollection(1).Formula = "=SERIES(,(Tabelle1!`$A`$2:`$A`$14,Tabelle1!`$A`$16:`$A`$19),(Tabelle1!`$C`$2:`$C`$14,Tabelle1!`$C`$16:`$C`$19),1)"

$chart3 = $ws.ChartObjects(3).Chart
$chart3.SeriesCollection(1).Formula = "=SERIES(,(Tabelle1!`$A`$2:`$A`$14,Tabelle1!`$A`$16:`$A`$19),(Tabelle1!`$D`$2:`$D`$14,Tabelle1!`$D`$16:`$D`$19),1)"

$wb.Save()
